# Generate Report for Handback
# Updates the localization-status workbook to reflect that the handback
# has been completed and is in sync with en-US: status text changes,
# handback timestamps are refreshed, and the stale "not latest" error
# detail is cleared out.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# The "Status" columns got wider (to fit the longer handed-back message)
# and the "Error Detail" columns got narrower (now that the cell is
# empty). Re-apply the closest achievable widths.
$statusColWidth = 29.166666666666668
$errorColWidth = 12.833333333333334

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $statusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $statusColWidth

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("K2").Value = "2016-09-03 06:52:38"
$wsZhCn.Range("P2").Value = ""
$wsZhCn.Columns.Item(3).ColumnWidth = $statusColWidth
$wsZhCn.Columns.Item(16).ColumnWidth = $errorColWidth

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-09-03 06:52:45"
$wsDeDe.Range("P2").Value = ""
$wsDeDe.Columns.Item(3).ColumnWidth = $statusColWidth
$wsDeDe.Columns.Item(16).ColumnWidth = $errorColWidth
